$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 6399
$ws.Range("J3").Value = 6795
$ws.Range("J4").Value = 1466
$ws.Range("J5").Value = 522
$ws.Range("J6").Value = 9006
$ws.Range("J7").Value = 24188

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 192
$ws.Range("J6").Value = 182
$ws.Range("J7").Value = 704
$ws.Range("J8").Value = 1528
$ws.Range("J10").Value = 179
$ws.Range("J12").Value = 49
$ws.Range("J14").Value = 126
$ws.Range("J15").Value = 280
$ws.Range("J19").Value = 708
$ws.Range("J20").Value = 504
$ws.Range("J25").Value = 118
$ws.Range("J29").Value = 1322
$ws.Range("J31").Value = 228
$ws.Range("J33").Value = 1089
$ws.Range("J35").Value = 32
$ws.Range("J36").Value = 329
$ws.Range("J37").Value = 742
$ws.Range("J41").Value = 163
$ws.Range("J42").Value = 1043
$ws.Range("J43").Value = 204
$ws.Range("J44").Value = 183
$ws.Range("J50").Value = 146
$ws.Range("J52").Value = 612
$ws.Range("J53").Value = 352
$ws.Range("J54").Value = 458
$ws.Range("J55").Value = 370
$ws.Range("J57").Value = 106
$ws.Range("J63").Value = 80
$ws.Range("J64").Value = 156
$ws.Range("J65").Value = 599
$ws.Range("J66").Value = 72
$ws.Range("J67").Value = 908
$ws.Range("J70").Value = 37
$ws.Range("J72").Value = 95
$ws.Range("J73").Value = 232
$ws.Range("J75").Value = 72
$ws.Range("J79").Value = 682
$ws.Range("J80").Value = 40
$ws.Range("J85").Value = 1004
$ws.Range("J87").Value = 79
$ws.Range("J91").Value = 274
$ws.Range("J94").Value = 252
$ws.Range("J95").Value = 349
$ws.Range("J97").Value = 217
$ws.Range("J100").Value = 45
$ws.Range("J101").Value = 24188

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J6").Value = 49
$ws.Range("J7").Value = 126

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J3").Value = 216
$ws.Range("J6").Value = 225
$ws.Range("J7").Value = 704

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 357
$ws.Range("J6").Value = 290
$ws.Range("J7").Value = 1004

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J6").Value = 261
$ws.Range("J7").Value = 612

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 63
$ws.Range("J6").Value = 233
$ws.Range("J7").Value = 352

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 413
$ws.Range("J3").Value = 458
$ws.Range("J5").Value = 39
$ws.Range("J6").Value = 536
$ws.Range("J7").Value = 1528

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 257
$ws.Range("J3").Value = 360
$ws.Range("J7").Value = 1089

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J3").Value = 126
$ws.Range("J7").Value = 349

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 252
$ws.Range("J6").Value = 219
$ws.Range("J7").Value = 742

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J6").Value = 215
$ws.Range("J7").Value = 599

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J2").Value = 82
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 228

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 339
$ws.Range("J6").Value = 251
$ws.Range("J7").Value = 908

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value = 215
$ws.Range("J7").Value = 458

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 466
$ws.Range("J7").Value = 1322

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J6").Value = 276
$ws.Range("J7").Value = 708

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J6").Value = 73
$ws.Range("J7").Value = 183

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J3").Value = 45
$ws.Range("J6").Value = 68
$ws.Range("J7").Value = 182

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J3").Value = 23
$ws.Range("J6").Value = 95
$ws.Range("J7").Value = 163

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 222
$ws.Range("J6").Value = 557
$ws.Range("J7").Value = 1043

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J6").Value = 100
$ws.Range("J7").Value = 179

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J3").Value = 74
$ws.Range("J6").Value = 205
$ws.Range("J7").Value = 370

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J6").Value = 68
$ws.Range("J7").Value = 274

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 189
$ws.Range("J3").Value = 229
$ws.Range("J6").Value = 203
$ws.Range("J7").Value = 682

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J2").Value = 42
$ws.Range("J7").Value = 156

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J6").Value = 136
$ws.Range("J7").Value = 504

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 107
$ws.Range("J7").Value = 329

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J4").Value = 18
$ws.Range("J7").Value = 252

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J2").Value = 53
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 118

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J3").Value = 64
$ws.Range("J7").Value = 280

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J2").Value = 38
$ws.Range("J3").Value = 37
$ws.Range("J7").Value = 146

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 72

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("J3").Value = 6
$ws.Range("J7").Value = 32

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J6").Value = 81
$ws.Range("J7").Value = 232

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J2").Value = 55
$ws.Range("J6").Value = 73
$ws.Range("J7").Value = 192

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J6").Value = 153
$ws.Range("J7").Value = 217

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("J2").Value = 15
$ws.Range("J6").Value = 5
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("J3").Value = 23
$ws.Range("J7").Value = 72

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J3").Value = 27
$ws.Range("J6").Value = 47
$ws.Range("J7").Value = 106

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J3").Value = 39
$ws.Range("J6").Value = 122
$ws.Range("J7").Value = 204

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 95

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("J4").Value = 7
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("J2").Value = 8
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 79
